$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A5").Value = "parisk"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "nan"
$ws.Range("D5").Value = "DIS"
$ws.Range("E5").Value = "WRI"
$ws.Range("F5").Value = "53dcf950-aee9-43ba-bb93-9e7c5cd5833d"
$ws.Range("G5").Value = "By5SY2gA-_annotated.xlsx"
$ws.Range("H5").Value = "For instance, what about averaging WordNet path-based distance metrics and distance in word embedding space (for word similarity), and other ways of applying the affect data to email tone prediction?"
$ws.Range("I5").Value = "Correct"
